$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Create the three new character styles (wdStyleTypeCharacter = 2)
# ---------------------------------------------------------------------
$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# ---------------------------------------------------------------------
# 2. Apply GaNStyle to every occurrence of the "V roku 2022 ..." run
# ---------------------------------------------------------------------
$rng = $d.Content
while ($rng.Find.Execute("V roku 2022 môžete pozorovať Súhvezdie Býk: 16. – 25. januára", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# ---------------------------------------------------------------------
# 3. Apply GaNParagraph to the "Stávate sa súčasťou ..." run
# ---------------------------------------------------------------------
$rng = $d.Content
if ($rng.Find.Execute("Stávate sa súčasťou celosvetovej kampane Globe at Night, ktorej cieľom je meranie svetelného znečistenia. Pozorovaním  Súhvezdie Býk na nočnej oblohe a porovnávaním skutočnej situácie s našimi mapkami sa nielenže dozviete, ako osvetlenie vo Vašom okolí prispieva k svetelnému znečisteniu, ale budete môcť porovnať úroveň svetelného znečistenia aj s inými lokalitami z celého sveta. Vaše pozorovanie tiež rozšíri online databázu dokumentujúcu viditeľnosť nočnej oblohy na našej planéte", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNParagraph"
}

# ---------------------------------------------------------------------
# 4. Apply GaNLinks to the "Mapky v tomto dokumente ..." run
# ---------------------------------------------------------------------
$rng = $d.Content
if ($rng.Find.Execute("Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNLinks"
}

Write-Host "Styles created and applied."
